# Auto-generated edit script: updates market-price-derived columns (H-N)
# for Halicarnassus_Profits sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2852.8333
$ws.Cells.Item(43, 10).Value = 2523.4
$ws.Cells.Item(43, 12).Value = 2523.4
$ws.Cells.Item(43, 14).Value = -2661.4

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 9079.714
$ws.Cells.Item(62, 9).Value = 6999
$ws.Cells.Item(62, 11).Value = 6999
$ws.Cells.Item(62, 13).Value = -6375

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 9079.714
$ws.Cells.Item(65, 9).Value = 6999
$ws.Cells.Item(65, 11).Value = 34995
$ws.Cells.Item(65, 13).Value = -31875

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 6984.143
$ws.Cells.Item(69, 9).Value = 6895
$ws.Cells.Item(69, 11).Value = 20685
$ws.Cells.Item(69, 13).Value = -19811

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 6984.143
$ws.Cells.Item(72, 9).Value = 6895
$ws.Cells.Item(72, 11).Value = 62055
$ws.Cells.Item(72, 13).Value = -57687

# ALC row 95
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 10124
$ws.Cells.Item(95, 10).Value = 10124
$ws.Cells.Item(95, 12).Value = 10124
$ws.Cells.Item(95, 14).Value = -15616

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 6430.273
$ws.Cells.Item(116, 9).Value = 5869.75
$ws.Cells.Item(116, 11).Value = 5869.75
$ws.Cells.Item(116, 13).Value = -2427.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 4970
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 4970
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 14910
$ws.Cells.Item(138, 13).Value = ""
$ws.Cells.Item(138, 14).Value = -25190

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3959.8262
$ws.Cells.Item(2, 9).Value = 3404.4
$ws.Cells.Item(2, 10).Value = 5001.25
$ws.Cells.Item(2, 11).Value = 3404.4
$ws.Cells.Item(2, 12).Value = 5001.25
$ws.Cells.Item(2, 13).Value = -3291.4
$ws.Cells.Item(2, 14).Value = -5227.25

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6628
$ws.Cells.Item(61, 9).Value = 6170.6665
$ws.Cells.Item(61, 10).Value = 8000
$ws.Cells.Item(61, 11).Value = 6170.6665
$ws.Cells.Item(61, 12).Value = 8000
$ws.Cells.Item(61, 13).Value = -5958.6665
$ws.Cells.Item(61, 14).Value = -8424

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 3219.1667
$ws.Cells.Item(88, 10).Value = 2250
$ws.Cells.Item(88, 12).Value = 2250
$ws.Cells.Item(88, 14).Value = -3062

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 3219.1667
$ws.Cells.Item(91, 10).Value = 2250
$ws.Cells.Item(91, 12).Value = 2250
$ws.Cells.Item(91, 14).Value = -5058

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1123.238
$ws.Cells.Item(97, 9).Value = 899.3333
$ws.Cells.Item(97, 10).Value = 2466.6667
$ws.Cells.Item(97, 11).Value = 899.3333
$ws.Cells.Item(97, 12).Value = 2466.6667
$ws.Cells.Item(97, 13).Value = -403.3333
$ws.Cells.Item(97, 14).Value = -3458.6667

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 4235.2
$ws.Cells.Item(102, 9).Value = 2794
$ws.Cells.Item(102, 11).Value = 2794
$ws.Cells.Item(102, 13).Value = -1172

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 3959.8262
$ws.Cells.Item(116, 9).Value = 3404.4
$ws.Cells.Item(116, 10).Value = 5001.25
$ws.Cells.Item(116, 11).Value = 3404.4
$ws.Cells.Item(116, 12).Value = 5001.25
$ws.Cells.Item(116, 13).Value = -1110.4
$ws.Cells.Item(116, 14).Value = -9589.25

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1750.8572
$ws.Cells.Item(122, 9).Value = 1564.25
$ws.Cells.Item(122, 10).Value = 1999.6666
$ws.Cells.Item(122, 11).Value = 4692.75
$ws.Cells.Item(122, 12).Value = 5998.9998
$ws.Cells.Item(122, 13).Value = -2242.75
$ws.Cells.Item(122, 14).Value = -10898.9998

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 6628
$ws.Cells.Item(136, 9).Value = 6170.6665
$ws.Cells.Item(136, 10).Value = 8000
$ws.Cells.Item(136, 11).Value = 18511.9995
$ws.Cells.Item(136, 12).Value = 24000
$ws.Cells.Item(136, 13).Value = -15961.9995
$ws.Cells.Item(136, 14).Value = -29100

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3959.8262
$ws.Cells.Item(3, 9).Value = 3404.4
$ws.Cells.Item(3, 10).Value = 5001.25
$ws.Cells.Item(3, 11).Value = 3404.4
$ws.Cells.Item(3, 12).Value = 5001.25
$ws.Cells.Item(3, 13).Value = -3290.4
$ws.Cells.Item(3, 14).Value = -5229.25

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 5972.364
$ws.Cells.Item(86, 9).Value = 3282.6667
$ws.Cells.Item(86, 11).Value = 3282.6667
$ws.Cells.Item(86, 13).Value = -2159.6667

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 5972.364
$ws.Cells.Item(89, 9).Value = 3282.6667
$ws.Cells.Item(89, 11).Value = 16413.3335
$ws.Cells.Item(89, 13).Value = -10797.3335

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 3245.7715
$ws.Cells.Item(107, 9).Value = 1547.56
$ws.Cells.Item(107, 10).Value = 7491.3
$ws.Cells.Item(107, 11).Value = 1547.56
$ws.Cells.Item(107, 12).Value = 7491.3
$ws.Cells.Item(107, 13).Value = 372.4400000000001
$ws.Cells.Item(107, 14).Value = -11331.3

# CRP row 5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 70
$ws.Cells.Item(5, 9).Value = 70
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 70
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 42
$ws.Cells.Item(5, 14).Value = ""

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 5424.75
$ws.Cells.Item(22, 10).Value = 5999.8
$ws.Cells.Item(22, 12).Value = 5999.8
$ws.Cells.Item(22, 14).Value = -6699.8

# CRP row 63
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(63, 8).Value = 40000
$ws.Cells.Item(63, 10).Value = 40000
$ws.Cells.Item(63, 12).Value = 40000
$ws.Cells.Item(63, 14).Value = -41372

# CRP row 66
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(66, 8).Value = 40000
$ws.Cells.Item(66, 10).Value = 40000
$ws.Cells.Item(66, 12).Value = 120000
$ws.Cells.Item(66, 14).Value = -126864

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 5855.1665
$ws.Cells.Item(99, 9).Value = 5027.3
$ws.Cells.Item(99, 10).Value = 9994.5
$ws.Cells.Item(99, 11).Value = 5027.3
$ws.Cells.Item(99, 12).Value = 9994.5
$ws.Cells.Item(99, 13).Value = -3529.3
$ws.Cells.Item(99, 14).Value = -12990.5

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 5855.1665
$ws.Cells.Item(126, 9).Value = 5027.3
$ws.Cells.Item(126, 10).Value = 9994.5
$ws.Cells.Item(126, 11).Value = 15081.9
$ws.Cells.Item(126, 12).Value = 29983.5
$ws.Cells.Item(126, 13).Value = -12611.9
$ws.Cells.Item(126, 14).Value = -34923.5

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3815.5
$ws.Cells.Item(134, 9).Value = 3214.3333
$ws.Cells.Item(134, 11).Value = 9642.999899999999
$ws.Cells.Item(134, 13).Value = -7107.999899999999

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 283.08694
$ws.Cells.Item(12, 9).Value = 155.85715
$ws.Cells.Item(12, 10).Value = 338.75
$ws.Cells.Item(12, 11).Value = 467.57145
$ws.Cells.Item(12, 12).Value = 1016.25
$ws.Cells.Item(12, 13).Value = -294.57145
$ws.Cells.Item(12, 14).Value = -1362.25

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2131.3333
$ws.Cells.Item(113, 10).Value = 2224
$ws.Cells.Item(113, 12).Value = 6672
$ws.Cells.Item(113, 14).Value = -11012

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 13).Value = ""
$ws.Cells.Item(129, 14).Value = ""

# GSM row 33
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 10500
$ws.Cells.Item(33, 10).Value = 10500
$ws.Cells.Item(33, 12).Value = 10500
$ws.Cells.Item(33, 14).Value = -11004

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5554.6113
$ws.Cells.Item(70, 9).Value = 3726.875
$ws.Cells.Item(70, 10).Value = 7016.8
$ws.Cells.Item(70, 11).Value = 3726.875
$ws.Cells.Item(70, 12).Value = 7016.8
$ws.Cells.Item(70, 13).Value = -3456.875
$ws.Cells.Item(70, 14).Value = -7556.8

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5554.6113
$ws.Cells.Item(73, 9).Value = 3726.875
$ws.Cells.Item(73, 10).Value = 7016.8
$ws.Cells.Item(73, 11).Value = 3726.875
$ws.Cells.Item(73, 12).Value = 7016.8
$ws.Cells.Item(73, 13).Value = -2790.875
$ws.Cells.Item(73, 14).Value = -8888.799999999999

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 5845.6665
$ws.Cells.Item(113, 9).Value = 3210.625
$ws.Cells.Item(113, 11).Value = 3210.625
$ws.Cells.Item(113, 13).Value = -1040.625

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2663.111
$ws.Cells.Item(122, 9).Value = 2023.875
$ws.Cells.Item(122, 11).Value = 6071.625
$ws.Cells.Item(122, 13).Value = -3621.625

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 128954.664
$ws.Cells.Item(132, 9).Value = 184365.33
$ws.Cells.Item(132, 10).Value = 18133.334
$ws.Cells.Item(132, 11).Value = 553095.99
$ws.Cells.Item(132, 12).Value = 54400.00199999999
$ws.Cells.Item(132, 13).Value = -550565.99
$ws.Cells.Item(132, 14).Value = -59460.00199999999

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5752.7646
$ws.Cells.Item(61, 9).Value = 5027
$ws.Cells.Item(61, 11).Value = 5027
$ws.Cells.Item(61, 13).Value = -4825

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 5752.7646
$ws.Cells.Item(113, 9).Value = 5027
$ws.Cells.Item(113, 11).Value = 5027
$ws.Cells.Item(113, 13).Value = -2857

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4732.5
$ws.Cells.Item(122, 9).Value = 4598.75
$ws.Cells.Item(122, 11).Value = 13796.25
$ws.Cells.Item(122, 13).Value = -11346.25

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 10181.125
$ws.Cells.Item(132, 9).Value = 9090.200000000001
$ws.Cells.Item(132, 10).Value = 11999.333
$ws.Cells.Item(132, 11).Value = 27270.6
$ws.Cells.Item(132, 12).Value = 35997.999
$ws.Cells.Item(132, 13).Value = -24740.6
$ws.Cells.Item(132, 14).Value = -41057.999

# WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 43291.668
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 43291.668
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 43291.668
$ws.Cells.Item(41, 13).Value = ""
$ws.Cells.Item(41, 14).Value = -44071.668

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 11966.667
$ws.Cells.Item(62, 9).Value = 11900
$ws.Cells.Item(62, 11).Value = 11900
$ws.Cells.Item(62, 13).Value = -11276

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 11966.667
$ws.Cells.Item(65, 9).Value = 11900
$ws.Cells.Item(65, 11).Value = 59500
$ws.Cells.Item(65, 13).Value = -56380

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 6790.8423
$ws.Cells.Item(136, 9).Value = 6190.125
$ws.Cells.Item(136, 10).Value = 9994.666999999999
$ws.Cells.Item(136, 11).Value = 18570.375
$ws.Cells.Item(136, 12).Value = 29984.001
$ws.Cells.Item(136, 13).Value = -16020.375
$ws.Cells.Item(136, 14).Value = -35084.001
